$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.06072676737724855
$ws.Range("C3").Value = 0.08943893963504915
$ws.Range("C4").Value = 0.1334090148888231
$ws.Range("C5").Value = 0.1676458654000764
$ws.Range("C6").Value = 0.1764930039625781
$ws.Range("C7").Value = 0.2169601972114539
$ws.Range("C8").Value = 0.2002337577397597
$ws.Range("C9").Value = 0.2005158181554662
$ws.Range("C10").Value = 0.2222691216053256
$ws.Range("C11").Value = 0.2852834459338497
$ws.Range("C12").Value = 0.3208952380823946
$ws.Range("C13").Value = 0.3930933767943192
$ws.Range("C14").Value = 0.4847235668517124
$ws.Range("C15").Value = 0.7094182127464319
$ws.Range("C16").Value = 0.08555299016961894
$ws.Range("C17").Value = 0.1145107702291279
$ws.Range("C18").Value = 0.1985598520837476
$ws.Range("C19").Value = 0.2624276034285558
$ws.Range("C20").Value = 0.2610714115743212
$ws.Range("C21").Value = 0.2792387315653621
$ws.Range("C22").Value = 0.2931935588784764
$ws.Range("C23").Value = 0.2705306887847534
$ws.Range("C24").Value = 0.3096114273698124
$ws.Range("C25").Value = 0.3806659143632873
$ws.Range("C26").Value = 0.468502703203439
$ws.Range("C27").Value = 0.5296433334641457
$ws.Range("C28").Value = 0.67243917541257
$ws.Range("C29").Value = 0.8949639471001851
$ws.Range("C30").Value = 0.04991793491499974
$ws.Range("C31").Value = 0.06203764819603896
$ws.Range("C32").Value = 0.1038756479498138
$ws.Range("C33").Value = 0.2136148753027318
$ws.Range("C34").Value = 0.1940082160738814
$ws.Range("C35").Value = 0.1815584906478788
$ws.Range("C36").Value = 0.1731002802120558
$ws.Range("C37").Value = 0.2406581184562085
$ws.Range("C38").Value = 0.2878955832071318
$ws.Range("C39").Value = 0.408218505485928
$ws.Range("C40").Value = 0.4714796346991122
$ws.Range("C41").Value = 0.5929192240979314
$ws.Range("C42").Value = 0.7951875975015328
$ws.Range("C43").Value = 1.115965900707341
